$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 entirely; row 3 shifts up to become the new row 2.
$ws.Rows("2").Delete()
